$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the HFR (004381180) and IVONE (004452476) rows (originally Excel rows 4-5)
$ws.Range("A4:C5").EntireRow.Delete()

# Remove the ANTONIO (004241147) and LARISSA (004846293) rows
# (originally rows 11-12, now rows 9-10 after the deletion above)
$ws.Range("A9:C10").EntireRow.Delete()

# Insert a new row right after VALMIR (004487140, now row 9) and populate
# it with the new WLADMIR (004388077) account
$ws.Range("A10:C10").EntireRow.Insert()

# Force the account number into the cell as text so the leading zeros
# survive (Excel would otherwise coerce the numeric-looking string to a
# number), then drop the temporary formatting so the cell matches the
# plain/unstyled look of its neighbours.
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "004388077"
$ws.Cells.Item(10, 1).ClearFormats()

$ws.Cells.Item(10, 2).Value = "WLADMIR"
$ws.Cells.Item(10, 3).Value = 5213.88
